# "Last bug fixes for monoclinic reindexing"
#
# Updates a batch of recomputed E-column statistics on the "Template"
# sheet (one brand-new cell, E21, plus corrections to 59 existing cells),
# mirrors that same newly-computed E21 value onto the matching row (E27)
# of the "Groups V0" sheet, and leaves the UI focused back on the
# "Template" sheet/cell E2 (with "Groups V0" parked at E27) the way the
# author's Excel session ended up.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Template")
$ws2 = $wb.Worksheets.Item("Groups V0")

# --- Template: corrected / newly-filled-in E-column values -----------------
$templateUpdates = [ordered]@{
    "E4"   = 4799
    "E6"   = 1163
    "E12"  = 688
    "E13"  = 17
    "E21"  = 2123   # brand-new cell (row previously had no E value)
    "E22"  = 9
    "E27"  = 772
    "E29"  = 1653
    "E30"  = 13
    "E31"  = 108
    "E32"  = 0
    "E33"  = 778
    "E34"  = 7
    "E35"  = 508
    "E38"  = 4893
    "E39"  = 1627
    "E40"  = 4892
    "E43"  = 4912
    "E44"  = 27
    "E45"  = 468
    "E49"  = 4871
    "E50"  = 0
    "E52"  = 2
    "E58"  = 21
    "E59"  = 545
    "E66"  = 3484
    "E67"  = 9
    "E70"  = 41
    "E72"  = 3156
    "E74"  = 1652
    "E75"  = 4
    "E78"  = 214
    "E79"  = 9
    "E83"  = 4891
    "E84"  = 2777
    "E85"  = 4901
    "E88"  = 4915
    "E89"  = 23
    "E90"  = 435
    "E92"  = 117
    "E94"  = 4813
    "E96"  = 4854
    "E97"  = 37
    "E100" = 5
    "E102" = 924
    "E103" = 39
    "E104" = 687
    "E112" = 10
    "E115" = 49
    "E117" = 847
    "E119" = 1645
    "E120" = 7
    "E121" = 124
    "E123" = 2871
    "E124" = 0
    "E125" = 1985
    "E128" = 4895
    "E129" = 3346
    "E133" = 4911
    "E134" = 128
    "E135" = 2934
}

foreach ($cellRef in $templateUpdates.Keys) {
    $ws1.Range($cellRef).Value = $templateUpdates[$cellRef]
}

# --- Groups V0: same newly-derived value lands on its matching row ---------
$ws2.Range("E27").Value = 2123

# --- Restore the window/selection state the workbook was left in -----------
# "Groups V0" was the active tab before; afterwards "Template" is active
# (with cell E2 selected) and "Groups V0" is left scrolled to/selecting E27.
$ws2.Activate() | Out-Null
$ws2.Range("E27").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("E2").Select() | Out-Null
